$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" value looks numeric need an explicit Text
# number format first, otherwise Excel auto-converts the literal into a
# real number instead of keeping it as the original text string.
$textCells = @("D5","D6","D9","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D27","D28","D29","D31","D34","D35","D37","D38","D39","D40","D41","D43","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.621.39"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.168.62"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "616.15"
$ws.Range("E5").Value = "  +3.19%  "
$ws.Range("D6").Value = "147.97"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.164.75"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "35.89"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "3.687.33"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "64.671.09"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "3.162.02"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "6.95"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "480.15"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "14.78"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "0.725"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "7.99"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("D24").Value = "13.79"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "84.52"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "8.59"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "26.59"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "53.03"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("D40").Value = "461.54"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "0.0401"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "8.43"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "2.851.50"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  +5.95%  "
$ws.Range("D48").Value = "26.69"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "120.21"
$ws.Range("E51").Value = "  +0.83%  "
